# Apply the "New crime data collected" weekly update to the 77th Precinct CompStat sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header strings: bulletin volume/week updated ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Crime statistics table (rows 14-30) ---
$ws.Range("M14").Value = -50
$ws.Range("M14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M14").Font.Bold = $false
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 20
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 28
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = 3.703703703703
$ws.Range("L16").Value = 154.545454545455
$ws.Range("M16").Value = -34.883720930232
$ws.Range("N16").Value = -91.194968553459
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 39
$ws.Range("K17").Value = 23.076923076923
$ws.Range("L17").Value = 54.838709677419
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = -65.217391304347
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = -46.153846153846
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -62.162162162162
$ws.Range("N18").Value = -92.820512820512
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 7.407407407407
$ws.Range("I19").Value = 59
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 18
$ws.Range("L19").Value = 51.282051282051
$ws.Range("M19").Value = 73.529411764705
$ws.Range("N19").Value = -20.270270270270
$ws.Range("C20").Value = 4
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Font.Bold = $false
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Font.Bold = $false
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 57.142857142857
$ws.Range("I20").Value = 22
$ws.Range("K20").Value = 29.411764705882
$ws.Range("L20").Value = 214.285714285714
$ws.Range("M20").Value = 29.411764705882
$ws.Range("N20").Value = -82.4
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 26.315789473684
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = 1.052631578947
$ws.Range("I21").Value = 178
$ws.Range("J21").Value = 164
$ws.Range("K21").Value = 8.536585365853
$ws.Range("L21").Value = 69.523809523809
$ws.Range("M21").Value = 7.228915662650
$ws.Range("N21").Value = -79.657142857142
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("F22").NumberFormat = "General"
$ws.Range("F22").Font.Bold = $false
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = -100
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -85.714285714285
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = 14.285714285714
$ws.Range("L23").Value = 128.571428571429
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 44.444444444444
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = 14.893617021276
$ws.Range("I24").Value = 92
$ws.Range("J24").Value = 94
$ws.Range("K24").Value = -2.127659574468
$ws.Range("L24").Value = 13.580246913580
$ws.Range("M24").Value = -10.679611650485
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -13.333333333333
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = 14.583333333333
$ws.Range("I25").Value = 94
$ws.Range("J25").Value = 75
$ws.Range("K25").Value = 25.333333333333
$ws.Range("L25").Value = 129.268292682927
$ws.Range("M25").Value = 4.444444444444
$ws.Range("D26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -85.714285714285
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -33.333333333333
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 4
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Font.Bold = $false
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Font.Bold = $false
$ws.Range("F27").Value = 10
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 11
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = -26.666666666666
$ws.Range("L27").Value = 57.142857142857
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Font.Bold = $false
$ws.Range("F28").Value = 1
$ws.Range("F28").NumberFormat = "#,##0"
$ws.Range("F28").Font.Bold = $false
$ws.Range("I28").Value = 1
$ws.Range("I28").NumberFormat = "#,##0"
$ws.Range("I28").Font.Bold = $false
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -85.714285714285
$ws.Range("N28").Value = -97.435897435897
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Font.Bold = $false
$ws.Range("F29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("F29").Font.Bold = $false
$ws.Range("I29").Value = 1
$ws.Range("I29").NumberFormat = "#,##0"
$ws.Range("I29").Font.Bold = $false
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -80
$ws.Range("N29").Value = -97.297297297297
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Font.Bold = $false
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Font.Bold = $false
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("G30").Font.Bold = $false
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H30").Font.Bold = $false
$ws.Range("J30").Value = 1
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("J30").Font.Bold = $false
$ws.Range("K30").Value = -100
$ws.Range("K30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K30").Font.Bold = $false
